# Locate the paragraph ending with " vers les 18h15, avant d'adapter mon code
# pour les autres types de données. " and append the new sentences after it,
# within the same paragraph (before the paragraph mark), matching the target
# diff exactly.

$d = $word.ActiveDocument

$marker = "vers les 18h15, avant"
$found = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$marker*") {
        $found = $p
        break
    }
}

if ($found -eq $null) {
    throw "Could not find target paragraph containing '$marker'"
}

$pRange = $found.Range
# Collapsed insertion point right before the paragraph mark (the very end of
# the paragraph's textual content).
$insertPoint = $d.Range($pRange.End - 1, $pRange.End - 1)

# First two runs are plain text runs appended in sequence, exactly as in the
# diff. A placeholder marker character ("~") is appended at the very end;
# it stands in for the remainder of the paragraph so that the following
# InsertXML call (needed to place the <w:proofErr/> markers) targets a
# non-collapsed range, which keeps the new content merged inside this same
# paragraph instead of being split into a new one.
$insertPoint.InsertAfter("Après une journée bien rempli j’ai finis d’adapter ")
$insertPoint.Collapse(0)
$insertPoint.InsertAfter("le ")
$insertPoint.Collapse(0)
$insertPoint.InsertAfter("~")
$insertPoint.Collapse(0)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$tailXml = '<w:p ' + $wNs + '>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>template</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">à gros coup de </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">if pour les 5 types supporté à l’heure sur mon site. </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Comme à chaque fois je vais donc push, actuellement il est 19h33 et le push sera sous le nom de « Adaptation de la page résultat en fonction du type de contenu ». </w:t></w:r>' +
    '</w:p>'

$tailRange = $d.Range($insertPoint.End - 1, $pRange.End - 1)
$tailRange.InsertXML($tailXml)

Write-Host "Edit applied"
